# Fruta / hortaliza, semanal
# Update rows with new Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M),
# Origen (O) and Precio $/Kg (P) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg)

    $ws.Cells.Item($Row, 4).Value  = $Fecha       # D: Fecha
    $ws.Cells.Item($Row, 10).Value = $Volumen     # J: Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin   # K: Precio minimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMax   # L: Precio maximo
    $ws.Cells.Item($Row, 13).Value = $PrecioProm  # M: Precio promedio ponderado
    $ws.Cells.Item($Row, 15).Value = $Origen      # O: Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg    # P: Precio $/Kg
}

Set-Row 2  44432 25 14000 14000 14000 "Provincia del Elquí" 467
Set-Row 3  44418 30 15000 15000 15000 "Provincia de Limarí" 500
Set-Row 4  44421 25 15000 16000 15400 "Provincia de Limarí" 513
Set-Row 5  44460 45 13000 13000 13000 "Provincia de Limarí" 433
Set-Row 6  44425 35 14000 14000 14000 "Provincia de Limarí" 467
Set-Row 7  44467 35 12000 12000 12000 "Provincia de Limarí" 400
Set-Row 8  44446 25 14000 14000 14000 "Provincia de Limarí" 467
Set-Row 11 44376 25 18000 18000 18000 "Provincia de Limarí" 600
Set-Row 13 44449 45 12000 12000 12000 "Provincia de Limarí" 400
Set-Row 14 44474 45 10000 10000 10000 "Provincia de Limarí" 333
